# Updated symbol list: refresh Price (column D) and Volume(1h) (column E)
# values for the crypto rows whose quotes changed.
# A leading apostrophe forces Excel to store the value as literal text,
# matching the original inlineStr cells (avoids numeric/percent coercion).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.73"
$ws.Range("E2").Value = "'-0.14%"
$ws.Range("D3").Value = "'32.65"
$ws.Range("E3").Value = "'1.49%"
$ws.Range("D4").Value = "'5.059"
$ws.Range("E4").Value = "'-1.23%"
$ws.Range("D5").Value = "'0.07722"
$ws.Range("E5").Value = "'-1.67%"
$ws.Range("D6").Value = "'2.106"
$ws.Range("E6").Value = "'-6.45%"
$ws.Range("D7").Value = "'7.903"
$ws.Range("E7").Value = "'1.15%"
$ws.Range("D8").Value = "'0.9236"
$ws.Range("E8").Value = "'-0.41%"
$ws.Range("E9").Value = "'0.34%"
$ws.Range("D10").Value = "'0.07984"
$ws.Range("E10").Value = "'2.16%"
$ws.Range("D11").Value = "'0.08572"
$ws.Range("E11").Value = "'-2.75%"
$ws.Range("D12").Value = "'0.03074"
$ws.Range("E12").Value = "'0.64%"
$ws.Range("D13").Value = "'0.09973"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("D14").Value = "'0.001520"
$ws.Range("E14").Value = "'0.77%"
$ws.Range("D15").Value = "'0.005649"
$ws.Range("E15").Value = "'-5.85%"
$ws.Range("D17").Value = "'3.475"
$ws.Range("E17").Value = "'0.33%"
$ws.Range("D18").Value = "'3.792"
$ws.Range("E18").Value = "'-0.19%"
$ws.Range("E19").Value = "'-4.16%"
$ws.Range("D20").Value = "'0.3338"
$ws.Range("E20").Value = "'2.00%"
$ws.Range("D21").Value = "'0.1318"
$ws.Range("E21").Value = "'-2.19%"
$ws.Range("D22").Value = "'4.393"
$ws.Range("E22").Value = "'3.50%"
$ws.Range("D23").Value = "'0.1973"
$ws.Range("E23").Value = "'9.81%"
$ws.Range("D24").Value = "'0.04541"
$ws.Range("E24").Value = "'-0.92%"
$ws.Range("D25").Value = "'0.001229"
$ws.Range("E25").Value = "'-1.94%"
$ws.Range("D26").Value = "'0.004152"
$ws.Range("E26").Value = "'-7.68%"
$ws.Range("D27").Value = "'0.0001249"
$ws.Range("E27").Value = "'-0.09%"
$ws.Range("D39").Value = "'0.01732"
$ws.Range("E39").Value = "'-3.25%"
$ws.Range("D40").Value = "'0.04712"
$ws.Range("E40").Value = "'-0.36%"
$ws.Range("D41").Value = "'0.007464"
$ws.Range("E41").Value = "'3.12%"
$ws.Range("E42").Value = "'-0.72%"
$ws.Range("D43").Value = "'0.002327"
$ws.Range("E43").Value = "'9.65%"
$ws.Range("D44").Value = "'0.01060"
$ws.Range("E44").Value = "'-4.36%"
$ws.Range("D45").Value = "'0.00006149"
$ws.Range("E45").Value = "'-0.91%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.28%"
$ws.Range("D47").Value = "'1.112"
$ws.Range("E47").Value = "'-0.47%"
$ws.Range("D48").Value = "'0.002997"
$ws.Range("E48").Value = "'-6.51%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.28%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.28%"
